$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 2.1
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 2.88
$ws.Range("L4").Value = 4.75
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("U4").Value = 2.2
$ws.Range("V4").Value = 1.62
$ws.Range("W4").Value = 5.5
$ws.Range("X4").Value = 8.5
$ws.Range("Z4").Value = 19
$ws.Range("AA4").Value = 21
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 19
$ws.Range("AF4").Value = 81
$ws.Range("AH4").Value = 8.5
$ws.Range("AI4").Value = 19
$ws.Range("AJ4").Value = 15
$ws.Range("AK4").Value = 41
$ws.Range("AL4").Value = 41
$ws.Range("AM4").Value = 51
$ws.Range("AN4").Value = 4
$ws.Range("AO4").Value = 12
$ws.Range("AQ4").Value = 41
$ws.Range("AS4").Value = 251
$ws.Range("AU4").Value = 9.5
$ws.Range("AV4").Value = 81
$ws.Range("AX4").Value = 5.5
$ws.Range("AY4").Value = 23
$ws.Range("AZ4").Value = 41
$ws.Range("BA4").Value = 81
$ws.Range("BB4").Value = 126
$ws.Range("BC4").Value = 351
$ws.Range("O5").Value = 1.19
$ws.Range("P5").Value = 4.33
$ws.Range("S5").Value = 1.32
$ws.Range("T5").Value = 3.25
$ws.Range("U5").Value = 1.62
$ws.Range("V5").Value = 2.21
$ws.Range("P6").Value = 4.06
$ws.Range("O7").Value = 1.18
$ws.Range("P7").Value = 4.21
$ws.Range("U7").Value = 1.67
$ws.Range("V7").Value = 2.15
$ws.Range("G8").Value = 1.65
$ws.Range("I8").Value = 5.5
$ws.Range("J8").Value = 2.27
$ws.Range("L8").Value = 5.8
$ws.Range("M8").Value = 1.09
$ws.Range("N8").Value = 7.25
$ws.Range("P8").Value = 2.35
$ws.Range("Q8").Value = 2.32
$ws.Range("U8").Value = 2.27
$ws.Range("W8").Value = 4.8
$ws.Range("X8").Value = 6.2
$ws.Range("Z8").Value = 11.5
$ws.Range("AA8").Value = 17
$ws.Range("AC8").Value = 6.7
$ws.Range("AH8").Value = 10.75
$ws.Range("AI8").Value = 30
$ws.Range("AJ8").Value = 19.5
$ws.Range("AK8").Value = 120
$ws.Range("AL8").Value = 80
$ws.Range("AM8").Value = 100
$ws.Range("AN8").Value = 3.2
$ws.Range("AO8").Value = 8.25
$ws.Range("AP8").Value = 23
$ws.Range("AQ8").Value = 30
$ws.Range("AT8").Value = 2.22
$ws.Range("AX8").Value = 6.7
$ws.Range("AY8").Value = 35
$ws.Range("AZ8").Value = 50
